$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.005583666666666666
$ws.Range("H2").Value = 0.016751
$ws.Range("I2").Value = 0.0005349789730684028
$ws.Range("J2").Value = 0.0005349789730684027
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06449866666666666
$ws.Range("N2").Value = 0.193496
$ws.Range("O2").Value = 0.004525829983623641
$ws.Range("P2").Value = 0.004525829983623642
$ws.Range("Q2").Value = 0.0003601390551111111
$ws.Range("R2").Value = 0.003241251496
$ws.Range("S2").Value = 0.000002421223876921162
$ws.Range("T2").Value = 0.000002421223876921162
# Row 3
$ws.Range("G3").Value = 0.005583666666666666
$ws.Range("H3").Value = 0.016751
$ws.Range("I3").Value = 0.0005349789730684028
$ws.Range("J3").Value = 0.0005349789730684027
$ws.Range("O3").Value = 0.745188142173877
$ws.Range("P3").Value = 0.7451881421738772
$ws.Range("Q3").Value = 0.05929770989488888
$ws.Range("R3").Value = 0.5336793890539999
$ws.Range("S3").Value = 0.0003986599870429317
$ws.Range("T3").Value = 0.0003986599870429317
# Row 4
$ws.Range("G4").Value = 0.005583666666666666
$ws.Range("H4").Value = 0.016751
$ws.Range("I4").Value = 0.0005349789730684028
$ws.Range("J4").Value = 0.0005349789730684027
$ws.Range("M4").Value = 3.566885000000001
$ws.Range("N4").Value = 10.700655
$ws.Range("O4").Value = 0.2502860278424993
$ws.Range("P4").Value = 0.2502860278424993
$ws.Range("Q4").Value = 0.01991629687833334
$ws.Range("R4").Value = 0.179246671905
$ws.Range("S4").Value = 0.0001338977621485499
$ws.Range("T4").Value = 0.0001338977621485499
# Row 5
$ws.Range("I5").Value = 0.7634962107066631
$ws.Range("J5").Value = 0.7634962107066631
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.06449866666666666
$ws.Range("N5").Value = 0.193496
$ws.Range("O5").Value = 0.004525829983623641
$ws.Range("P5").Value = 0.004525829983623642
$ws.Range("Q5").Value = 0.5139731050133334
$ws.Range("R5").Value = 4.62575794512
$ws.Range("S5").Value = 0.003455454042799249
$ws.Range("T5").Value = 0.00345545404279925
# Row 6
$ws.Range("I6").Value = 0.7634962107066631
$ws.Range("J6").Value = 0.7634962107066631
$ws.Range("O6").Value = 0.745188142173877
$ws.Range("P6").Value = 0.7451881421738772
$ws.Range("S6").Value = 0.5689483228132932
$ws.Range("T6").Value = 0.5689483228132933
# Row 7
$ws.Range("I7").Value = 0.7634962107066631
$ws.Range("J7").Value = 0.7634962107066631
$ws.Range("M7").Value = 3.566885000000001
$ws.Range("N7").Value = 10.700655
$ws.Range("O7").Value = 0.2502860278424993
$ws.Range("P7").Value = 0.2502860278424993
$ws.Range("Q7").Value = 28.42357917490001
$ws.Range("R7").Value = 255.8122125741
$ws.Range("S7").Value = 0.1910924338505706
$ws.Range("T7").Value = 0.1910924338505706
# Row 8
$ws.Range("G8").Value = 2.462846666666667
$ws.Range("H8").Value = 7.388540000000001
$ws.Range("I8").Value = 0.2359688103202685
$ws.Range("J8").Value = 0.2359688103202685
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.06449866666666666
$ws.Range("N8").Value = 0.193496
$ws.Range("O8").Value = 0.004525829983623641
$ws.Range("P8").Value = 0.004525829983623642
$ws.Range("Q8").Value = 0.1588503262044444
$ws.Range("R8").Value = 1.42965293584
$ws.Range("S8").Value = 0.001067954716947471
$ws.Range("T8").Value = 0.001067954716947471
# Row 9
$ws.Range("G9").Value = 2.462846666666667
$ws.Range("H9").Value = 7.388540000000001
$ws.Range("I9").Value = 0.2359688103202685
$ws.Range("J9").Value = 0.2359688103202685
$ws.Range("O9").Value = 0.745188142173877
$ws.Range("P9").Value = 0.7451881421738772
$ws.Range("Q9").Value = 26.15506545679556
$ws.Range("R9").Value = 235.39558911116
$ws.Range("S9").Value = 0.1758411593735409
$ws.Range("T9").Value = 0.1758411593735409
# Row 10
$ws.Range("G10").Value = 2.462846666666667
$ws.Range("H10").Value = 7.388540000000001
$ws.Range("I10").Value = 0.2359688103202685
$ws.Range("J10").Value = 0.2359688103202685
$ws.Range("M10").Value = 3.566885000000001
$ws.Range("N10").Value = 10.700655
$ws.Range("O10").Value = 0.2502860278424993
$ws.Range("P10").Value = 0.2502860278424993
$ws.Range("Q10").Value = 8.784690832633336
$ws.Range("R10").Value = 79.06221749370002
$ws.Range("S10").Value = 0.05905969622978015
$ws.Range("T10").Value = 0.05905969622978016
